# Generate Report for Handoff
#
# For the six files that just reached "Ready for handoff" (rows 7, 8, 9,
# 10, 13, 14 in each per-language sheet), the handoff-xliff-generation
# step now records:
#   - a Priority ("ht") in column E of the zh-cn / de-de sheets
#   - refreshed "Latest HO Xliff Generate Date" / "Latest Handoff Datetime"
#     timestamps in the Overview sheet (col G) and the zh-cn / de-de sheets
#     (col H)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$rows = @(7, 8, 9, 10, 13, 14)

foreach ($r in $rows) {
    # Priority column (E) on the per-language sheets: blank -> "ht"
    $zhcn.Range("E$r").Value = "ht"
    $dede.Range("E$r").Value = "ht"

    # Latest HO Xliff Generate Date (Overview!G) / Latest Handoff Datetime
    # (de-de!H) shared the same timestamp; bump it forward.
    $overview.Range("G$r").Value = "2016-08-26 18:22:37"
    $dede.Range("H$r").Value = "2016-08-26 18:22:37"

    # Latest Handoff Datetime on zh-cn!H had its own timestamp.
    $zhcn.Range("H$r").Value = "2016-08-26 18:22:32"
}
